# Update the "想去人数" (want-to-go count) figures for two events.
# These values appear duplicated on both the "展览" sheet and the
# "全部类型" sheet, so both need to be updated to stay consistent.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 5331
    $ws.Range("F4").Value = 920
}
